# Update cryptos list with latest price/volume figures (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: several Price (D) values are plain numeric-looking strings. Excel would
# normally auto-convert these to numbers on assignment; the source column is
# text, so we force a Text format while assigning, then drop back to the
# Normal style so no stray formatting is left on the cell.

$ws.Range("D2").Value = "45.429.14"
$ws.Range("E2").Value = "  +7.09%  "
$ws.Range("D3").Value = "2.379.01"
$ws.Range("E3").Value = "  +4.18%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "113.28"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +9.95%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "317.25"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.47%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.633"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.04%  "
$ws.Range("E8").Value = "  -0.25%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.625"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.19%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.63"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +10.63%  "
$ws.Range("E11").Value = "  +3.43%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.66"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +5.84%  "
$ws.Range("B13").Value = "Polygon"
$ws.Range("C13").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.02"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.07%  "
$ws.Range("B14").Value = "TRON"
$ws.Range("C14").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.109"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.50%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.86"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.92%  "
$ws.Range("D16").Value = "2.741.95"
$ws.Range("E16").Value = "  +4.38%  "
$ws.Range("D17").Value = "2.379.90"
$ws.Range("E17").Value = "  +4.55%  "
$ws.Range("D18").Value = "45.318.00"
$ws.Range("E18").Value = "  +6.33%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.63"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.68%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000107"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.45%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.36"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.76%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "74.83"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.74%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.53"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.46%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "269.25"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.62%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.35"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +8.96%  "
$ws.Range("E26").Value = "  -0.53%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.61"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +8.35%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.31"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +6.03%  "
$ws.Range("E29").Value = "  -0.11%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "39.14"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +9.58%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.91"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.82%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0966"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +14.03%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "171.88"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.50%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.96"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +15.84%  "
$ws.Range("E35").Value = "  +2.54%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.120"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +8.52%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.93"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +10.04%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.07"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +12.57%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.01"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +11.21%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0365"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +6.06%  "
$ws.Range("E41").Value = "  +10.68%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "104.98"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.23%  "
$ws.Range("E43").Value = "  +6.77%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "71.45"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.45%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.22"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +10.39%  "
$ws.Range("E46").Value = "  -0.34%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.80"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +12.80%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "116.49"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +6.01%  "
$ws.Range("E49").Value = "  +18.47%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.34"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +8.08%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "79.09"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.55%  "
